# Insert a new weekly price record as row 19 on the "Orégano" sheet.
# All the existing rows 19:34 shift down to 20:35 (EntireRow.Insert takes
# care of that, along with the sheet's used-range / dimension growing from
# R34 to R35), and the new row 19 is populated with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 19:34 down one slot to make room for the new weekly record.
$ws.Rows(19).Insert()

# Fill in the new row 19 with this week's values (same shape as its
# neighbours: Mercado/Región/Categoría metadata repeated, only the date and
# price columns change).
$ws.Cells.Item(19, 1).Value = 9
$ws.Cells.Item(19, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(19, 3).Value = "Metropolitana"
$ws.Cells.Item(19, 4).Value = 44489
$ws.Cells.Item(19, 5).Value = 13
$ws.Cells.Item(19, 6).Value = 100112029
$ws.Cells.Item(19, 7).Value = "Orégano"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 16
$ws.Cells.Item(19, 11).Value = 9000
$ws.Cells.Item(19, 12).Value = 10000
$ws.Cells.Item(19, 13).Value = 9500
$ws.Cells.Item(19, 14).Value = "$/docena de atados"
$ws.Cells.Item(19, 15).Value = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value = 3167
$ws.Cells.Item(19, 17).Value = 3
$ws.Cells.Item(19, 18).Value = "Hortaliza"
